$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (preserve existing formatting) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "45"

$weekCell = $ws.Range("C9")
$weekCell.Characters(48, 9).Text = "11/10/2024"
$weekCell.Characters(27, 10).Text = "11/4/2024"

# --- Update crime statistics grid (rows 16-31) ---
$ws.Range("C16").Value = 2
$ws.Range("I14").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 5
$ws.Range("K14").Copy($ws.Range("E16"))
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 65
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = -35.643564356435
$ws.Range("L16").Value = -18.75
$ws.Range("M16").Value = 18.181818181818
$ws.Range("N16").Value = -88.203266787658
$ws.Range("C17").Value = 1
$ws.Range("C14").Copy($ws.Range("D17"))
$ws.Range("E14").Copy($ws.Range("E17"))
$ws.Range("F17").Value = 5
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 25
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 96
$ws.Range("K17").Value = 3.125
$ws.Range("L17").Value = -3.883495145631
$ws.Range("M17").Value = 102.040816326531
$ws.Range("N17").Value = -28.260869565217
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 12.5
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 107
$ws.Range("K18").Value = -26.168224299065
$ws.Range("L18").Value = -50.931677018633
$ws.Range("M18").Value = -3.658536585365
$ws.Range("N18").Value = -92.518939393939
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 7.317073170731
$ws.Range("I19").Value = 512
$ws.Range("J19").Value = 600
$ws.Range("K19").Value = -14.666666666666
$ws.Range("L19").Value = -11.111111111111
$ws.Range("M19").Value = -20.124804992199
$ws.Range("N19").Value = -73.957273652085
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 2
$ws.Range("C14").Copy($ws.Range("G20"))
$ws.Range("E14").Copy($ws.Range("H20"))
$ws.Range("I20").Value = 23
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = -52.083333333333
$ws.Range("L20").Value = -60.344827586206
$ws.Range("M20").Value = -20.689655172413
$ws.Range("N20").Value = -95.855855855855
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -18.181818181818
$ws.Range("F21").Value = 65
$ws.Range("G21").Value = 61
$ws.Range("H21").Value = 6.557377049180
$ws.Range("I21").Value = 788
$ws.Range("J21").Value = 959
$ws.Range("K21").Value = -17.831074035453
$ws.Range("L21").Value = -20.564516129032
$ws.Range("M21").Value = -8.584686774941
$ws.Range("N21").Value = -81.580177653108
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 2
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 3
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -25
$ws.Range("I22").Value = 33
$ws.Range("J22").Value = 35
$ws.Range("K22").Value = -5.714285714285
$ws.Range("L22").Value = 26.923076923076
$ws.Range("M22").Value = 10
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -62.068965517241
$ws.Range("F24").Value = 53
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -50.925925925925
$ws.Range("I24").Value = 812
$ws.Range("J24").Value = 920
$ws.Range("K24").Value = -11.739130434782
$ws.Range("L24").Value = -26.449275362318
$ws.Range("M24").Value = 49.539594843462
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -52.631578947368
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 84
$ws.Range("H25").Value = -59.523809523809
$ws.Range("I25").Value = 627
$ws.Range("J25").Value = 759
$ws.Range("K25").Value = -17.391304347826
$ws.Range("L25").Value = -31.773667029379
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 600
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 201
$ws.Range("J26").Value = 215
$ws.Range("K26").Value = -6.511627906976
$ws.Range("L26").Value = -2.427184466019
$ws.Range("M26").Value = 1.005025125628
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
$ws.Range("L28").Value = -2.083333333333
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
